# Update the dSF column (F) values to reflect the re-pulled data.
# These rows had their F value diverge from the original copy of column E (dS0).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = 3
    6  = 0
    10 = -7
    15 = -1
    18 = -5
    19 = 0
    23 = -1
    28 = -4
    30 = -6
    32 = -1
    36 = -10
    45 = 12
    47 = -6
    52 = -4
    53 = -3
    60 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
